$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '61.262.46'
$ws.Range('E2').Value = '  -1.64%  '
$ws.Range('D3').Value = '2.983.80'
$ws.Range('E3').Value = '  -0.45%  '
$ws.Range('E4').Value = '  -0.10%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '599.00'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +3.14%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '143.76'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -1.52%  '
$ws.Range('E7').Value = '  -0.03%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.519'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.60%  '
$ws.Range('D9').Value = '2.981.81'
$ws.Range('E9').Value = '  -0.60%  '
$ws.Range('B10').Value = 'Dogecoin'
$ws.Range('C10').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.146'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -1.32%  '
$ws.Range('B11').Value = 'Toncoin'
$ws.Range('C11').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '6.05'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +7.08%  '
$ws.Range('E12').Value = '  +3.08%  '
$ws.Range('E13').Value = '  -0.08%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '34.29'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -0.71%  '
$ws.Range('E15').Value = '  +2.38%  '
$ws.Range('D16').Value = '3.478.25'
$ws.Range('E16').Value = '  -0.49%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '6.93'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -1.90%  '
$ws.Range('D18').Value = '61.232.93'
$ws.Range('E18').Value = '  -1.73%  '
$ws.Range('D19').Value = '2.980.33'
$ws.Range('E19').Value = '  -0.56%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '449.13'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -1.49%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '14.20'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +2.48%  '
$ws.Range('E22').Value = '  +0.66%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.33'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.51%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '81.92'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +2.45%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.19'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -3.56%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '10.47'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +4.85%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '11.93'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -2.65%  '
$ws.Range('E28').Value = '  +0.08%  '
$ws.Range('E29').Value = '  +3.21%  '
$ws.Range('E30').Value = '  -0.03%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '7.15'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.42%  '
$ws.Range('E32').Value = '  -1.75%  '
$ws.Range('E33').Value = '  +0.90%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.109'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +1.88%  '
$ws.Range('D35').Value = '0.0₃0823'
$ws.Range('E35').Value = '  +5.01%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.01'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.46%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '5.78'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.99%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '50.31'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.48%  '
$ws.Range('E39').Value = '  -2.32%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '9.10'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +0.71%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.123'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +10.55%  '
$ws.Range('E42').Value = '  -1.38%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '395.86'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -3.46%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '39.94'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +4.62%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0350'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.10%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.268'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -2.85%  '
$ws.Range('D47').Value = '2.687.33'
$ws.Range('E47').Value = '  -2.97%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '131.29'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E49').Value = '  +0.09%  '
$ws.Range('E50').Value = '  -0.65%  '
$ws.Range('E51').Value = '  +0.57%  '
